# "Editing the Docs for Main table"
# Mark the Main table's Docs column (D) with "V" for the Phaser_Output (row 8)
# and Phases (row 9) entries, matching the rest of the Main_Table rows above.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D8").Value = "V"
$ws.Range("D9").Value = "V"

# Re-select the full table (as last left by the author) with D10 as the
# anchor cell that was clicked before selecting the whole range.
$ws.Range("D10").Activate()
$ws.Range("A1:D14").Select()
